# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# before column N, pushing the old N/O/P ("Late" / heading / "Outstanding")
# columns one to the right (O/P/Q). The sheet also becomes the active tab
# with a new single-cell selection, and the former active sheet
# ("Edit Repayment Schedule") loses its tab-selected state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column); this shifts the
# existing N:P columns (and their widths/styles) one column to the right.
$ws.Columns("N:N").Insert()

# The new column N ends up with the sheet's default width; give it the same
# rendered width as column M (11 characters) like the authored workbook.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and set its selection.
$ws.Activate() | Out-Null
$ws.Range("K18").Select() | Out-Null
